$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 91, shifting all existing data rows (91-230) down by one
# (to rows 92-231). This mirrors the weekly price update commit, which adds a new
# price observation for the latest date at the top of the historical series.
$ws.Rows.Item(91).Insert()

$newRow = 91
$ws.Cells.Item($newRow, 1).Value2  = 5
$ws.Cells.Item($newRow, 2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item($newRow, 3).Value2  = "Maule"
$ws.Cells.Item($newRow, 4).Value2  = 44665
$ws.Cells.Item($newRow, 5).Value2  = 7
$ws.Cells.Item($newRow, 6).Value2  = 100112009
$ws.Cells.Item($newRow, 7).Value2  = "Acelga"
$ws.Cells.Item($newRow, 8).Value2  = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value2  = "Primera"
$ws.Cells.Item($newRow, 10).Value2 = 500
$ws.Cells.Item($newRow, 11).Value2 = 3500
$ws.Cells.Item($newRow, 12).Value2 = 3500
$ws.Cells.Item($newRow, 13).Value2 = 3500
$ws.Cells.Item($newRow, 14).Value2 = "`$/docena de atados (4 kilos)"
$ws.Cells.Item($newRow, 15).Value2 = "Región del Maule"
$ws.Cells.Item($newRow, 16).Value2 = 875
$ws.Cells.Item($newRow, 17).Value2 = 4
$ws.Cells.Item($newRow, 18).Value2 = "Hortaliza"

# Make sure the date cell keeps the same date/time number format as the rest of
# column D.
$ws.Cells.Item($newRow, 4).NumberFormat = $ws.Cells.Item($newRow + 1, 4).NumberFormat
